$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (all target cells are plain text
# in the sheet, incl. numeric-looking price strings like "30.015.74" / "0.9997"
# and percentage strings like "  -0.12%  ").
$updates = [ordered]@{
    'D2' = '30.015.74'
    'D3' = '1.883.31'
    'E3' = '  -0.12%  '
    'D4' = '0.9997'
    'E4' = '  -0.04%  '
    'D5' = '0.7354'
    'E5' = '  -1.81%  '
    'D6' = '241.94'
    'E6' = '  -0.05%  '
    'D7' = '0.9996'
    'E7' = '  -0.10%  '
    'D8' = '0.3154'
    'E8' = '  +1.16%  '
    'D9' = '0.07159'
    'E9' = '  +0.93%  '
    'E10' = '  -1.99%  '
    'D11' = '0.08310'
    'E11' = '  -2.18%  '
    'D12' = '0.7554'
    'E12' = '  -0.34%  '
    'D13' = '5.392'
    'E13' = '  +0.74%  '
    'D14' = '1.849.19'
    'E14' = '  -2.28%  '
    'E15' = '  -0.37%  '
    'D16' = '6.138'
    'E16' = '  +0.28%  '
    'D17' = '30.009.77'
    'E17' = '  +1.14%  '
    'D18' = '248.50'
    'E18' = '  +2.46%  '
    'D19' = '13.54'
    'E19' = '  -0.98%  '
    'D20' = '0.000007838'
    'E20' = '  +0.21%  '
    'D21' = '2.141.00'
    'E21' = '  +0.06%  '
    'D22' = '0.9987'
    'E22' = '  -0.05%  '
    'D23' = '0.9993'
    'E23' = '  -0.08%  '
    'D24' = '7.882'
    'E24' = '  -0.73%  '
    'D25' = '0.1563'
    'E25' = '  -1.42%  '
    'E26' = '  -0.84%  '
    'D27' = '163.88'
    'E27' = '  +0.64%  '
    'E28' = '  -0.01%  '
    'E29' = '  +0.92%  '
    'D30' = '1.473'
    'E30' = '  +0.16%  '
    'D31' = '4.560'
    'E31' = '  +1.40%  '
    'D32' = '1.533'
    'E32' = '  +0.05%  '
    'D33' = '4.187'
    'E33' = '  +0.56%  '
    'D34' = '0.05312'
    'E34' = '  -1.72%  '
    'E35' = '  +0.73%  '
    'D36' = '0.7684'
    'E36' = '  +2.45%  '
    'D37' = '0.9973'
    'E37' = '  -0.74%  '
    'D38' = '2.722'
    'E38' = '  +0.55%  '
    'D39' = '0.01956'
    'E39' = '  +0.84%  '
    'D40' = '2.757'
    'E40' = '  -0.42%  '
    'D41' = '0.4575'
    'E41' = '  +2.64%  '
    'B42' = 'Maker'
    'C42' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D42' = '1.090.24'
    'E42' = '  -0.78%  '
    'B43' = 'FraxShare'
    'C43' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D43' = '6.031'
    'E43' = '  -0.81%  '
    'B44' = 'TrustWalletToken'
    'C44' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D44' = '0.8793'
    'E44' = '  +2.22%  '
    'D45' = '72.17'
    'E45' = '  -0.12%  '
    'D46' = '104.21'
    'E46' = '  +1.89%  '
    'D47' = '1.000'
    'E47' = '  -0.03%  '
    'E48' = '  -0.01%  '
    'D49' = '7.544'
    'E49' = '  -2.20%  '
    'D50' = '9.537'
    'E50' = '  -1.73%  '
    'D51' = '2.021.64'
    'E51' = '  -0.42%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "0.9997", "1.000")
    # are not auto-coerced into numbers by the COM layer.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default (unstyled) look so no stray number-format/style
    # is left behind on cells that originally had none.
    $cell.Style = "Normal"
}
